# Fixed update to excel issue
# The forecast was regenerated a week later than before, so every weekly
# row on "Forecast Comparison" shifts its Week_Start_Date forward by one
# week and gets fresh MyForecast / Amazon Mean / P70 / P80 / P90 numbers.
# The "Summary" sheet's derived statistics are recomputed to match.

$wb = $excel.ActiveWorkbook

$wsForecast = $wb.Worksheets.Item("Forecast Comparison")
$wsSummary  = $wb.Worksheets.Item("Summary")

# Week_Start_Date (col B) text values must stay literal text (not get
# auto-converted to Excel date serials), so every "B<row>" write is
# prefixed with a quote-prefix apostrophe, matching what a user typing
# a literal date string into a General-formatted cell would get.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
}

# row -> (Week_Start_Date, MyForecast, AmazonMean, AmazonP70, AmazonP80, AmazonP90)
$rows = @(
    @{ Row = 2;  Date = "2025-02-02"; D = 94;  E = 72;  F = 88;  G = 113; H = 155 },
    @{ Row = 3;  Date = "2025-02-09"; D = 82;  E = 63;  F = 77;  G = 99;  H = 134 },
    @{ Row = 4;  Date = "2025-02-16"; D = 101; E = 78;  F = 95;  G = 121; H = 163 },
    @{ Row = 5;  Date = "2025-02-23"; D = 114; E = 88;  F = 108; G = 134; H = 176 },
    @{ Row = 6;  Date = "2025-03-02"; D = 118; E = 91;  F = 111; G = 140; H = 187 },
    @{ Row = 7;  Date = "2025-03-09"; D = 110; E = 85;  F = 104; G = 133; H = 182 },
    @{ Row = 8;  Date = "2025-03-16"; D = 118; E = 91;  F = 111; G = 143; H = 195 },
    @{ Row = 9;  Date = "2025-03-23"; D = 117; E = 92;  F = 112; G = 145; H = 201 },
    @{ Row = 10; Date = "2025-03-30"; D = 108; E = 83;  F = 101; G = 131; H = 181 },
    @{ Row = 11; Date = "2025-04-06"; D = 112; E = 86;  F = 105; G = 137; H = 191 },
    @{ Row = 12; Date = "2025-04-13"; D = 114; E = 89;  F = 108; G = 142; H = 197 },
    @{ Row = 13; Date = "2025-04-20"; D = 110; E = 88;  F = 107; G = 139; H = 194 },
    @{ Row = 14; Date = "2025-04-27"; D = 103; E = 82;  F = 100; G = 130; H = 181 },
    @{ Row = 15; Date = "2025-05-04"; D = 97;  E = 80;  F = 98;  G = 128; H = 180 },
    @{ Row = 16; Date = "2025-05-11"; D = 91;  E = 84;  F = 102; G = 133; H = 184 },
    @{ Row = 17; Date = "2025-05-18"; D = 86;  E = 75;  F = 92;  G = 120; H = 168 }
)

foreach ($r in $rows) {
    $row = $r.Row
    Set-TextValue $wsForecast.Cells.Item($row, 2) $r.Date   # B: Week_Start_Date
    $wsForecast.Cells.Item($row, 4).Value = $r.D            # D: MyForecast
    $wsForecast.Cells.Item($row, 5).Value = $r.E            # E: Amazon Mean Forecast
    $wsForecast.Cells.Item($row, 6).Value = $r.F            # F: Amazon P70 Forecast
    $wsForecast.Cells.Item($row, 7).Value = $r.G            # G: Amazon P80 Forecast
    $wsForecast.Cells.Item($row, 8).Value = $r.H            # H: Amazon P90 Forecast
}

# Summary sheet recomputed statistics (all stored as literal text, same
# as every other "Value" cell on this sheet).
Set-TextValue $wsSummary.Range("B2")  "2022-12-25 to 2025-01-26"
Set-TextValue $wsSummary.Range("B4")  "309"
Set-TextValue $wsSummary.Range("B5")  "133"
Set-TextValue $wsSummary.Range("B6")  "135"
Set-TextValue $wsSummary.Range("B7")  "82"
Set-TextValue $wsSummary.Range("B8")  "14800 units"
Set-TextValue $wsSummary.Range("B9")  "1676"
Set-TextValue $wsSummary.Range("B10") "855"
Set-TextValue $wsSummary.Range("B11") "391"
Set-TextValue $wsSummary.Range("B12") "118"
Set-TextValue $wsSummary.Range("B13") "2025-03-02"
Set-TextValue $wsSummary.Range("B14") "82"
Set-TextValue $wsSummary.Range("B15") "2025-02-09"
